# CS304 Project.xlsx -- "Add files via upload"
#
# Re-enacts the author's edit pass over the workbook:
#  1. On the "product" sheet, bump the unit price (column C) for every
#     product row (2-73) up by a small, irregular amount.
#  2. Leave the cursor sitting on E12 of "product" (the last cell touched).
#  3. Switch over to the "supplierTerm" sheet and leave the cursor on I14 --
#     that sheet ends up being the active tab when the file was saved.

$wb = $excel.ActiveWorkbook

# --- 1) Update product prices (column C) on the "product" sheet ---
$wsProduct = $wb.Worksheets.Item("product")
$wsProduct.Activate() | Out-Null

$wsProduct.Range("C2").Value = 80.5
$wsProduct.Range("C3").Value = 55.5
$wsProduct.Range("C4").Value = 60.5
$wsProduct.Range("C5").Value = 61.5
$wsProduct.Range("C6").Value = 92.5
$wsProduct.Range("C7").Value = 81.5
$wsProduct.Range("C8").Value = 91
$wsProduct.Range("C9").Value = 70.5
$wsProduct.Range("C10").Value = 80.5
$wsProduct.Range("C11").Value = 101.5
$wsProduct.Range("C12").Value = 91.5
$wsProduct.Range("C13").Value = 47.5
$wsProduct.Range("C14").Value = 57.5
$wsProduct.Range("C15").Value = 82.5
$wsProduct.Range("C16").Value = 51.5
$wsProduct.Range("C17").Value = 143
$wsProduct.Range("C18").Value = 81.5
$wsProduct.Range("C19").Value = 82.5
$wsProduct.Range("C20").Value = 92.5
$wsProduct.Range("C21").Value = 41.5
$wsProduct.Range("C22").Value = 42.5
$wsProduct.Range("C23").Value = 61.5
$wsProduct.Range("C24").Value = 91.5
$wsProduct.Range("C25").Value = 26.5
$wsProduct.Range("C26").Value = 97.5
$wsProduct.Range("C27").Value = 60.5
$wsProduct.Range("C28").Value = 60.5
$wsProduct.Range("C29").Value = 25.88
$wsProduct.Range("C30").Value = 72.5
$wsProduct.Range("C31").Value = 77.5
$wsProduct.Range("C32").Value = 81
$wsProduct.Range("C33").Value = 142
$wsProduct.Range("C34").Value = 161
$wsProduct.Range("C35").Value = 97
$wsProduct.Range("C36").Value = 111
$wsProduct.Range("C37").Value = 40
$wsProduct.Range("C38").Value = 212
$wsProduct.Range("C39").Value = 150
$wsProduct.Range("C40").Value = 169
$wsProduct.Range("C41").Value = 183
$wsProduct.Range("C42").Value = 81
$wsProduct.Range("C43").Value = 601.99
$wsProduct.Range("C44").Value = 42.99
$wsProduct.Range("C45").Value = 15.99
$wsProduct.Range("C46").Value = 42.99
$wsProduct.Range("C47").Value = 35.99
$wsProduct.Range("C48").Value = 36.99
$wsProduct.Range("C49").Value = 26.99
$wsProduct.Range("C50").Value = 35.99
$wsProduct.Range("C51").Value = 12.99
$wsProduct.Range("C52").Value = 36.99
$wsProduct.Range("C53").Value = 35.99
$wsProduct.Range("C54").Value = 26.99
$wsProduct.Range("C55").Value = 71.989999999999995
$wsProduct.Range("C56").Value = 82.99
$wsProduct.Range("C57").Value = 26.99
$wsProduct.Range("C58").Value = 15.99
$wsProduct.Range("C59").Value = 62.99
$wsProduct.Range("C60").Value = 36.99
$wsProduct.Range("C61").Value = 20.99
$wsProduct.Range("C62").Value = 35.99
$wsProduct.Range("C63").Value = 20.99
$wsProduct.Range("C64").Value = 20.99
$wsProduct.Range("C65").Value = 15.99
$wsProduct.Range("C66").Value = 61.99
$wsProduct.Range("C67").Value = 50.99
$wsProduct.Range("C68").Value = 101.98
$wsProduct.Range("C69").Value = 18.989999999999998
$wsProduct.Range("C70").Value = 15.99
$wsProduct.Range("C71").Value = 50.99
$wsProduct.Range("C72").Value = 58.8
$wsProduct.Range("C73").Value = 51.2

# --- 2) Leave the selection where the author last clicked ---
$wsProduct.Range("E12").Select() | Out-Null

# --- 3) Navigate to "supplierTerm", which becomes the active tab on save ---
$wsSupplierTerm = $wb.Worksheets.Item("supplierTerm")
$wsSupplierTerm.Activate() | Out-Null
$wsSupplierTerm.Range("I14").Select() | Out-Null

